# daily auto push: 2026-02-12 14:15 UTC
#
# A new sample was recorded for 2026/02/12 (木/Thursday) at hour "22",
# appended after the existing 2026/02/12 rows (..., 13, 18) and before the
# 2026/12/29 block. Insert one new row at row 802 (shifting every
# subsequent row down by one, 802-843 -> 803-844) and populate it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 802..843 down to 803..844, leaving a blank row 802.
$ws.Rows.Item(802).Insert()

# A802 looks like a date ("2026/02/12") but must stay a literal text
# string (matching every other cell in column A). A leading apostrophe
# forces text interpretation instead of date-serial conversion; then
# re-apply the plain "Normal" style so no stray number-format sticks to
# the cell (every sibling row cell is unstyled).
$ws.Cells.Item(802, 1).Value = "'2026/02/12"
$ws.Cells.Item(802, 1).Style = "Normal"

$ws.Cells.Item(802, 2).Value = "木"
$ws.Cells.Item(802, 3).Value = 22
$ws.Cells.Item(802, 4).Value = 201
